$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, so Excel
# keeps them as literal text (matching the inlineStr content in the source).
$textCells = @(
    "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20",
    "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30",
    "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42",
    "D43", "D44", "D47", "D48", "D49", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range("D2").Value = "62.860.14"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "3.062.85"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "536.54"
$ws.Range("E5").Value = "  -3.85%  "
$ws.Range("D6").Value = "133.69"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.059.43"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "6.22"
$ws.Range("E11").Value = "  -7.79%  "
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "0.0000224"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Value = "34.40"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "3.555.11"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "62.795.97"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "3.063.48"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "6.63"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "483.90"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").Value = "13.30"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "7.15"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "79.42"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "12.17"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("D28").Value = "8.16"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "25.96"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "1.87"
$ws.Range("E31").Value = "  -8.70%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "2.38"
$ws.Range("E33").Value = "  -6.95%  "
$ws.Range("D34").Value = "56.35"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "477.41"
$ws.Range("E37").Value = "  -12.04%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0396"
$ws.Range("E38").Value = "  -4.70%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.086.83"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "0.0797"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "0.116"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "8.10"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "0.254"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D46").Value = "0.0₃0547"
$ws.Range("E46").Value = "  +8.82%  "
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("D48").Value = "120.78"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "24.67"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  +5.37%  "
